# Bugfix replay: "Se anadio en el i -- en el while de la funcion rangos ya
# que tenia ese bug" -- the membership-degree table ("rangos" output) was
# recomputed after fixing a missing loop-counter increment in a while loop,
# which shifted the GV frio / GV fresco membership values for several rows.
#
# The sheet stores these membership numbers as literal text (e.g. "0.8"),
# so each Range.Value assignment below uses a leading apostrophe to force
# the value to be stored as text rather than being auto-coerced to a
# number -- mirroring the original inline-string cell type for these
# columns. The Style is then reset to "Normal" right after so the
# text-marker (quote-prefix) formatting doesn't linger as a visible change.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 ("frio" bucket for value 0): GV frio / GV MAX 0.8 -> 0.2
$ws.Range("C2").Value = "'0.2"
$ws.Range("C2").Style = "Normal"
$ws.Range("G2").Value = "'0.2"
$ws.Range("G2").Style = "Normal"

# Row 5 ("fresco" bucket for value 3): recomputed membership numbers
$ws.Range("B5").Value = 9.890000000000001
$ws.Range("C5").Value = "'0.021999999999999888"
$ws.Range("C5").Style = "Normal"
$ws.Range("D5").Value = "'0.3150000000000001"
$ws.Range("D5").Style = "Normal"
$ws.Range("G5").Value = "'0.3150000000000001"
$ws.Range("G5").Style = "Normal"

# Row 6 ("frio" bucket for value 4): GV frio / GV MAX 0.2 -> 0.0
$ws.Range("C6").Value = "'0.0"
$ws.Range("C6").Style = "Normal"
$ws.Range("G6").Value = "'0.0"
$ws.Range("G6").Style = "Normal"

# Row 7 ("frio" bucket for value 5): GV frio / GV MAX 0.13333333333333333 -> 0.0
$ws.Range("C7").Value = "'0.0"
$ws.Range("C7").Style = "Normal"
$ws.Range("G7").Value = "'0.0"
$ws.Range("G7").Style = "Normal"

# Row 8 ("frio" bucket for value 6): GV frio / GV MAX 0.47600000000000003 -> 0.0
$ws.Range("C8").Value = "'0.0"
$ws.Range("C8").Style = "Normal"
$ws.Range("G8").Value = "'0.0"
$ws.Range("G8").Style = "Normal"
